$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("categoria") contains values such as "Compromisso" which must be
# renamed to "Convencao" throughout the used data range.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
$rangeE = $ws.Range("E1:E$lastRow")

$rangeE.Replace("Compromisso", "Convencao", 1, 1, $false, $false, $false, $false)

$wb.Save()
